# Update reaction time distributions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: Hydrogenation reaction time - Lower bound 1.9 -> 6.1, and swap highlight style to match row 43 (green fill)
$ws.Range("E35").Value = 6.1
$ws.Range("A43:K43").Copy() | Out-Null
$ws.Range("A35:K35").PasteSpecial(-4122) | Out-Null

# Row 38: change H38 formula from =G38 to =E38
$ws.Range("H38").Formula = "=E38"

# Row 41: Etherification & hydrolysis reaction time - Lower bound 1.9 -> 6.1, swap style to match row 42 (no fill)
$ws.Range("E41").Value = 6.1
$ws.Range("A42:K42").Copy() | Out-Null
$ws.Range("A41:K41").PasteSpecial(-4122) | Out-Null

# Row 46: Ring-opening & hydrolysis reaction time - Lower bound 1.9 -> 6.1
$ws.Range("E46").Value = 6.1

# Update sheet view: scroll position and selection
$excel.Goto($ws.Range("A24"), $true) | Out-Null
$ws.Range("A34:XFD54").Select() | Out-Null
